$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.546.90'
$ws.Range('E2').Value = '  +1.08%  '

# Row 3
$ws.Range('D3').Value = '1.794.63'
$ws.Range('E3').Value = '  -0.63%  '

# Row 4
$ws.Range('E4').Value = '  +0.33%  '

# Row 5
$ws.Range('D5').Value = '327.64'
$ws.Range('E5').Value = '  -3.11%  '

# Row 6
$ws.Range('E6').Value = '  +0.54%  '

# Row 7
$ws.Range('D7').Value = '0.4402'
$ws.Range('E7').Value = '  -3.40%  '

# Row 8
$ws.Range('D8').Value = '0.3725'
$ws.Range('E8').Value = '  +5.25%  '

# Row 9
$ws.Range('D9').Value = '45.72'
$ws.Range('E9').Value = '  +0.34%  '

# Row 10
$ws.Range('D10').Value = '0.07561'
$ws.Range('E10').Value = '  -0.62%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.130'
$ws.Range('E11').Value = '  -1.77%  '

# Row 12
$ws.Range('D12').Value = '22.58'
$ws.Range('E12').Value = '  -0.87%  '

# Row 13
$ws.Range('D13').Value = '1.006'
$ws.Range('E13').Value = '  +0.36%  '

# Row 14
$ws.Range('D14').Value = '6.199'
$ws.Range('E14').Value = '  -0.91%  '

# Row 15
$ws.Range('E15').Value = '  +3.08%  '

# Row 16
$ws.Range('D16').Value = '1.796.98'
$ws.Range('E16').Value = '  -0.51%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001084'
$ws.Range('E17').Value = '  -0.77%  '

# Row 18
$ws.Range('D18').Value = '0.06692'
$ws.Range('E18').Value = '  +0.25%  '

# Row 19
$ws.Range('D19').Value = '80.41'
$ws.Range('E19').Value = '  -1.79%  '

# Row 20
$ws.Range('E20').Value = '  +0.34%  '

# Row 21
$ws.Range('D21').Value = '17.45'
$ws.Range('E21').Value = '  +1.30%  '

# Row 22
$ws.Range('D22').Value = '6.212'
$ws.Range('E22').Value = '  -2.83%  '

# Row 23
$ws.Range('D23').Value = '28.565.01'
$ws.Range('E23').Value = '  +1.04%  '

# Row 24
$ws.Range('D24').Value = '11.65'
$ws.Range('E24').Value = '  -2.76%  '

# Row 25
$ws.Range('D25').Value = '2.437'
$ws.Range('E25').Value = '  +1.48%  '

# Row 26
$ws.Range('D26').Value = '20.38'
$ws.Range('E26').Value = '  -1.77%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.90'
$ws.Range('E27').Value = '  -1.65%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.330'
$ws.Range('E28').Value = '  -3.98%  '

# Row 29
$ws.Range('D29').Value = '2.004.16'
$ws.Range('E29').Value = '  -0.46%  '

# Row 30
$ws.Range('D30').Value = '1.308'
$ws.Range('E30').Value = '  +0.87%  '

# Row 31
$ws.Range('D31').Value = '130.47'
$ws.Range('E31').Value = '  -2.17%  '

# Row 32
$ws.Range('D32').Value = '3.982'
$ws.Range('E32').Value = '  -2.18%  '

# Row 33
$ws.Range('D33').Value = '5.766'
$ws.Range('E33').Value = '  -3.15%  '

# Row 34
$ws.Range('D34').Value = '0.09227'
$ws.Range('E34').Value = '  -2.69%  '

# Row 35
$ws.Range('D35').Value = '0.2236'
$ws.Range('E35').Value = '  +2.89%  '

# Row 36
$ws.Range('E36').Value = '  -1.27%  '

# Row 37
$ws.Range('D37').Value = '0.06237'
$ws.Range('E37').Value = '  -0.28%  '

# Row 38
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '5.176'
$ws.Range('E38').Value = '  -0.43%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02300'
$ws.Range('E39').Value = '  -3.66%  '

# Row 40
$ws.Range('D40').Value = '0.6549'
$ws.Range('E40').Value = '  -2.95%  '

# Row 41
$ws.Range('D41').Value = '1.196'
$ws.Range('E41').Value = '  -1.74%  '

# Row 42
$ws.Range('D42').Value = '1.428'
$ws.Range('E42').Value = '  -3.97%  '

# Row 43
$ws.Range('D43').Value = '7.963'
$ws.Range('E43').Value = '  -2.62%  '

# Row 44
$ws.Range('D44').Value = '1.003'
$ws.Range('E44').Value = '  +0.52%  '

# Row 45
$ws.Range('D45').Value = '13.89'
$ws.Range('E45').Value = '  -1.06%  '

# Row 46
$ws.Range('D46').Value = '0.6049'
$ws.Range('E46').Value = '  -1.15%  '

# Row 47
$ws.Range('D47').Value = '3.799'
$ws.Range('E47').Value = '  -1.68%  '

# Row 48
$ws.Range('D48').Value = '127.08'
$ws.Range('E48').Value = '  -2.11%  '

# Row 49
$ws.Range('D49').Value = '2.001'
$ws.Range('E49').Value = '  -1.87%  '

# Row 50
$ws.Range('D50').Value = '0.07005'
$ws.Range('E50').Value = '  -1.53%  '

# Row 51
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '1.134'
$ws.Range('E51').Value = '  -2.84%  '
